$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 27, shifting rows 27:84 down to 28:85
$ws.Rows("27:27").Insert()

# Fill in the new row 27 with data
$ws.Range("A27").Value = 3
$ws.Range("B27").Value = "Femacal de La Calera"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44498
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 100112026
$ws.Range("G27").Value = "Haba"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = 8000
$ws.Range("N27").Value = "$/malla 25 kilos"
$ws.Range("O27").Value = "Provincia de Quillota"
$ws.Range("P27").Value = 320
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = "Hortaliza"
